$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of trade data (row 8), mirroring the structure of the
# existing rows (3-7): Date, Profitable, Principle, Start Principle,
# BuyPrice, SellPrice, IsShortSell, Price Change %, Strong trade

# Copy formatting (styles) from row 7 down into row 8 first, so the new
# row's cells pick up the same style indexes (e.g. the date style used in
# columns A and G) without minting new number-format entries.
$ws.Range("A7:I7").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)

$ws.Range("A8").Value = 42650.366863425923

$ws.Range("B8").Value = $true

$ws.Range("C8").Value = 10345.49

$ws.Range("D8").Value = 10268.48

$ws.Range("E8").Value = 308.29998799999998

$ws.Range("F8").Value = 305.98998999999998

$ws.Range("G8").Value = $true

$ws.Range("H8").Value = -0.75

$ws.Range("I8").Value = $true
